$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $val) {
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.ClearFormats()
}

$ws.Range("D2").Value = '28.151.72'
$ws.Range("E2").Value = '  +0.43%  '

$ws.Range("D3").Value = '1.798.64'
$ws.Range("E3").Value = '  +2.61%  '

Set-TextValue $ws.Range("D4") '1.004'
$ws.Range("E4").Value = '  -0.14%  '

Set-TextValue $ws.Range("D5") '337.53'
$ws.Range("E5").Value = '  +0.19%  '

$ws.Range("E6").Value = '  +0.33%  '

Set-TextValue $ws.Range("D7") '0.4623'
$ws.Range("E7").Value = '  +22.75%  '

Set-TextValue $ws.Range("D8") '0.3696'
$ws.Range("E8").Value = '  +10.30%  '

Set-TextValue $ws.Range("D9") '45.23'
$ws.Range("E9").Value = '  -0.14%  '

Set-TextValue $ws.Range("D10") '0.07666'
$ws.Range("E10").Value = '  +6.61%  '

Set-TextValue $ws.Range("D11") '1.147'
$ws.Range("E11").Value = '  +2.73%  '

$ws.Range("B12").Value = 'Solana'
$ws.Range("C12").Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
Set-TextValue $ws.Range("D12") '22.52'
$ws.Range("E12").Value = '  +0.06%  '

$ws.Range("B13").Value = 'BinanceUSD'
$ws.Range("C13").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
Set-TextValue $ws.Range("D13") '1.003'
$ws.Range("E13").Value = '  +0.00%  '

Set-TextValue $ws.Range("D14") '6.356'
$ws.Range("E14").Value = '  +3.42%  '

Set-TextValue $ws.Range("D15") '7.394'
$ws.Range("E15").Value = '  +3.66%  '

$ws.Range("D16").Value = '1.795.56'
$ws.Range("E16").Value = '  +2.16%  '

Set-TextValue $ws.Range("D17") '0.00001094'
$ws.Range("E17").Value = '  +3.65%  '

Set-TextValue $ws.Range("D18") '0.06723'
$ws.Range("E18").Value = '  +2.07%  '

Set-TextValue $ws.Range("D19") '82.65'
$ws.Range("E19").Value = '  +2.76%  '

$ws.Range("E20").Value = '  +0.11%  '

Set-TextValue $ws.Range("D21") '17.43'
$ws.Range("E21").Value = '  +3.41%  '

Set-TextValue $ws.Range("D22") '6.410'
$ws.Range("E22").Value = '  +2.81%  '

$ws.Range("D23").Value = '28.139.12'
$ws.Range("E23").Value = '  +0.21%  '

Set-TextValue $ws.Range("D24") '11.90'
$ws.Range("E24").Value = '  +2.22%  '

Set-TextValue $ws.Range("D25") '2.417'
$ws.Range("E25").Value = '  +1.27%  '

Set-TextValue $ws.Range("D26") '20.78'
$ws.Range("E26").Value = '  +4.84%  '

Set-TextValue $ws.Range("D27") '2.386'
$ws.Range("E27").Value = '  +3.01%  '

Set-TextValue $ws.Range("D28") '152.16'
$ws.Range("E28").Value = '  -0.60%  '

$ws.Range("D29").Value = '2.001.70'
$ws.Range("E29").Value = '  +2.26%  '

Set-TextValue $ws.Range("D30") '133.77'
$ws.Range("E30").Value = '  +1.42%  '

Set-TextValue $ws.Range("D31") '1.262'
$ws.Range("E31").Value = '  +1.73%  '

Set-TextValue $ws.Range("D32") '4.053'
$ws.Range("E32").Value = '  +0.94%  '

Set-TextValue $ws.Range("D33") '0.09625'
$ws.Range("E33").Value = '  +10.68%  '

Set-TextValue $ws.Range("D34") '5.901'
$ws.Range("E34").Value = '  +2.14%  '

Set-TextValue $ws.Range("D35") '0.02388'
$ws.Range("E35").Value = '  +2.93%  '

Set-TextValue $ws.Range("D36") '0.2229'
$ws.Range("E36").Value = '  +5.86%  '

Set-TextValue $ws.Range("D37") '12.17'
$ws.Range("E37").Value = '  -0.09%  '

$ws.Range("B38").Value = 'InternetComputer(DFINITY)'
$ws.Range("C38").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
Set-TextValue $ws.Range("D38") '5.275'
$ws.Range("E38").Value = '  +2.37%  '

Set-TextValue $ws.Range("D39") '0.06343'
$ws.Range("E39").Value = '  +2.55%  '

$ws.Range("B40").Value = 'TheSandbox'
$ws.Range("C40").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
Set-TextValue $ws.Range("D40") '0.6727'
$ws.Range("E40").Value = '  +0.93%  '

Set-TextValue $ws.Range("D41") '1.507'
$ws.Range("E41").Value = '  +4.26%  '

$ws.Range("E42").Value = '  +1.56%  '

Set-TextValue $ws.Range("D43") '8.089'
$ws.Range("E43").Value = '  +1.00%  '

$ws.Range("E44").Value = '  +3.12%  '

$ws.Range("E45").Value = '  +0.21%  '

Set-TextValue $ws.Range("D46") '0.6172'
$ws.Range("E46").Value = '  +2.19%  '

$ws.Range("E47").Value = '  +0.38%  '

Set-TextValue $ws.Range("D48") '130.52'
$ws.Range("E48").Value = '  +1.65%  '

Set-TextValue $ws.Range("D49") '2.059'
$ws.Range("E49").Value = '  +2.23%  '

Set-TextValue $ws.Range("D50") '1.186'
$ws.Range("E50").Value = '  +1.40%  '

Set-TextValue $ws.Range("D51") '0.07128'
$ws.Range("E51").Value = '  -0.44%  '
